$wb = $excel.ActiveWorkbook

# hunk @1189 - sheet ALC row 11
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 48864.438
$ws.Range("I11").Value = 48864.438
$ws.Range("K11").Value = 48864.438
$ws.Range("M11").Value = -48724.438

# hunk @2729 - sheet ALC row 42
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 130.28572
$ws.Range("I42").Value = 130.28572
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 390.85716
$ws.Range("L42").Value = 0
$ws.Range("N42").Value = -160.85716

# hunk @3167 - sheet ALC row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3090.0688
$ws.Range("I51").Value = 3005.625
$ws.Range("J51").Value = 3277.7222
$ws.Range("K51").Value = 3005.625
$ws.Range("L51").Value = 3277.7222
$ws.Range("M51").Value = -2521.625
$ws.Range("N51").Value = -4245.7222

# hunk @3810 - sheet ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 7199.8
$ws.Range("I64").Value = 3999
$ws.Range("K64").Value = 3999
$ws.Range("M64").Value = -3751

# hunk @3960 - sheet ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 7199.8
$ws.Range("I67").Value = 3999
$ws.Range("K67").Value = 3999
$ws.Range("M67").Value = -3141

# hunk @5188 - sheet ALC row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 239.2381
$ws.Range("I92").Value = 239.2381
$ws.Range("K92").Value = 239.2381
$ws.Range("M92").Value = 1008.7619

# hunk @5883 - sheet ALC row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 17642.715
$ws.Range("I106").Value = 7999
$ws.Range("J106").Value = 19250
$ws.Range("K106").Value = 7999
$ws.Range("L106").Value = 19250
$ws.Range("M106").Value = -7368
$ws.Range("N106").Value = -20512

# hunk @6235 - sheet ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4098

# hunk @7169 - sheet ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 12479.091
$ws.Range("I132").Value = 17984.54
$ws.Range("K132").Value = 53953.62
$ws.Range("M132").Value = -51423.62

# hunk @9168 - sheet ARM row 31
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 3281.8333
$ws.Range("I31").Value = 3281.8333
$ws.Range("K31").Value = 3281.8333
$ws.Range("M31").Value = -2987.8333

# hunk @9217 - sheet ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 180466.1
$ws.Range("I32").Value = 186566.3
$ws.Range("J32").Value = 12710.5
$ws.Range("K32").Value = 186566.3
$ws.Range("L32").Value = 12710.5
$ws.Range("M32").Value = -186279.3
$ws.Range("N32").Value = -13284.5

# hunk @9860 - sheet ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3378.8
$ws.Range("I45").Value = 2298
$ws.Range("J45").Value = 5000
$ws.Range("K45").Value = 2298
$ws.Range("L45").Value = 5000
$ws.Range("M45").Value = -1921
$ws.Range("N45").Value = -5754

# hunk @10632 - sheet ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 11666.111
$ws.Range("I61").Value = 9285
$ws.Range("K61").Value = 9285
$ws.Range("M61").Value = -9073

# hunk @13564 - sheet ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2899.6667
$ws.Range("I122").Value = 2830.1
$ws.Range("J122").Value = 3247.5
$ws.Range("K122").Value = 8490.299999999999
$ws.Range("L122").Value = 9742.5
$ws.Range("M122").Value = -6040.299999999999
$ws.Range("N122").Value = -14642.5

# hunk @13616 - sheet ARM row 123
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 35000
$ws.Range("J123").Value = 35000
$ws.Range("L123").Value = 35000
$ws.Range("N123").Value = -44800

# hunk @14042 - sheet ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 760412.2
$ws.Range("I132").Value = 895727.7
$ws.Range("K132").Value = 2687183.1
$ws.Range("M132").Value = -2684653.1

# hunk @14241 - sheet ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 11666.111
$ws.Range("I136").Value = 9285
$ws.Range("K136").Value = 27855
$ws.Range("M136").Value = -25305

# hunk @15508 - sheet BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 94033.37
$ws.Range("I20").Value = 146824.28
$ws.Range("K20").Value = 146824.28
$ws.Range("M20").Value = -146577.28

# hunk @19386 - sheet BSM row 100
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 30000
$ws.Range("J100").Value = 30000
$ws.Range("L100").Value = 30000
$ws.Range("N100").Value = -32164

# hunk @19628 - sheet BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2618.5527
$ws.Range("I105").Value = 1823.3462
$ws.Range("K105").Value = 1823.3462
$ws.Range("M105").Value = -76.34619999999995

# hunk @21007 - sheet BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4128.2856
$ws.Range("I134").Value = 1936.9302
$ws.Range("K134").Value = 5810.7906
$ws.Range("M134").Value = -3275.7906

# hunk @21298 - sheet BSM row 140
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 79999
$ws.Range("J140").Value = 79999
$ws.Range("L140").Value = 79999
$ws.Range("N140").Value = -90359

# hunk @26231 - sheet CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4500
$ws.Range("J99").Value = 3600
$ws.Range("L99").Value = 3600
$ws.Range("N99").Value = -6596

# hunk @27346 - sheet CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2361.2273
$ws.Range("I122").Value = 1861.9286
$ws.Range("K122").Value = 5585.7858
$ws.Range("M122").Value = -3135.7858

# hunk @27539 - sheet CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 4500
$ws.Range("J126").Value = 3600
$ws.Range("L126").Value = 10800
$ws.Range("N126").Value = -15740

# hunk @27922 - sheet CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1646.0588
$ws.Range("I134").Value = 1682.2
$ws.Range("J134").Value = 1375
$ws.Range("K134").Value = 5046.6
$ws.Range("L134").Value = 4125
$ws.Range("M134").Value = -2511.6
$ws.Range("N134").Value = -9195

# hunk @30936 - sheet CUL row 52
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 3700
$ws.Range("J52").Value = 3700
$ws.Range("L52").Value = 11100
$ws.Range("N52").Value = -11632

# hunk @31487 - sheet CUL row 63
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("N63").Value = 0

# hunk @31637 - sheet CUL row 66
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("N66").Value = 0

# hunk @32369 - sheet CUL row 81
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 3375
$ws.Range("J81").Value = 5000
$ws.Range("L81").Value = 15000
$ws.Range("N81").Value = -17246

# hunk @32522 - sheet CUL row 84
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 3375
$ws.Range("J84").Value = 5000
$ws.Range("L84").Value = 45000
$ws.Range("N84").Value = -56232

# hunk @33174 - sheet CUL row 97
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 369.83334
$ws.Range("I97").Value = 146.66667
$ws.Range("J97").Value = 593
$ws.Range("K97").Value = 440.00001
$ws.Range("L97").Value = 1779
$ws.Range("M97").Value = 55.99998999999997
$ws.Range("N97").Value = -2771

# hunk @33379 - sheet CUL row 101
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H101").Value = 12510938
$ws.Range("J101").Value = 14296786
$ws.Range("L101").Value = 42890358
$ws.Range("N101").Value = -42895226

# hunk @34793 - sheet CUL row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1114729.1
$ws.Range("J129").Value = 3488.8333
$ws.Range("L129").Value = 10466.4999
$ws.Range("N129").Value = -20466.4999

# hunk @36656 - sheet GSM row 24
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 105099.6
$ws.Range("J24").Value = 5666.222
$ws.Range("L24").Value = 5666.222
$ws.Range("N24").Value = -6012.222

# hunk @38874 - sheet GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 87770.22
$ws.Range("J70").Value = 86748.125
$ws.Range("L70").Value = 86748.125
$ws.Range("N70").Value = -87288.125

# hunk @39018 - sheet GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 87770.22
$ws.Range("J73").Value = 86748.125
$ws.Range("L73").Value = 86748.125
$ws.Range("N73").Value = -88620.125

# hunk @40921 - sheet GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3532.5557
$ws.Range("I113").Value = 2550
$ws.Range("J113").Value = 5497.6665
$ws.Range("K113").Value = 2550
$ws.Range("L113").Value = 5497.6665
$ws.Range("M113").Value = -380
$ws.Range("N113").Value = -9837.666499999999

# hunk @41840 - sheet GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 8941.529
$ws.Range("I132").Value = 12403.317
$ws.Range("J132").Value = 3684.7407
$ws.Range("K132").Value = 37209.951
$ws.Range("L132").Value = 11054.2221
$ws.Range("M132").Value = -34679.951
$ws.Range("N132").Value = -16114.2221

# hunk @42663 - sheet LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5168.091
$ws.Range("J7").Value = 6797.8
$ws.Range("L7").Value = 6797.8
$ws.Range("N7").Value = -7021.8

# hunk @48196 - sheet LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5550.8696
$ws.Range("I122").Value = 4578.3335
$ws.Range("K122").Value = 13735.0005
$ws.Range("M122").Value = -11285.0005

# hunk @48392 - sheet LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5168.091
$ws.Range("J126").Value = 6797.8
$ws.Range("L126").Value = 20393.4
$ws.Range("N126").Value = -25333.4

# hunk @48683 - sheet LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1615036.6
$ws.Range("I132").Value = 2175897.2
$ws.Range("J132").Value = 2562.25
$ws.Range("K132").Value = 6527691.600000001
$ws.Range("L132").Value = 7686.75
$ws.Range("M132").Value = -6525161.600000001
$ws.Range("N132").Value = -12746.75

# hunk @51402 - sheet WVR row 46
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 60000
$ws.Range("J46").Value = 60000
$ws.Range("L46").Value = 60000
$ws.Range("N46").Value = -60462

# hunk @54331 - sheet WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1568.0286
$ws.Range("I107").Value = 982.0741
$ws.Range("K107").Value = 2946.2223
$ws.Range("M107").Value = -1026.2223

# hunk @55057 - sheet WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 29961.3
$ws.Range("I122").Value = 1742.4814
$ws.Range("J122").Value = 88569.62
$ws.Range("K122").Value = 5227.4442
$ws.Range("L122").Value = 265708.86
$ws.Range("M122").Value = -2777.4442
$ws.Range("N122").Value = -270608.86

# hunk @55645 - sheet WVR row 134
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 60000
$ws.Range("J134").Value = 60000
$ws.Range("L134").Value = 180000
$ws.Range("N134").Value = -185070

# hunk @55694 - sheet WVR row 135
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 79999
$ws.Range("J135").Value = 79999
$ws.Range("L135").Value = 79999
$ws.Range("N135").Value = -90139

# hunk @55740 - sheet WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1132.625
$ws.Range("I136").Value = 1198.25
$ws.Range("J136").Value = 935.75
$ws.Range("K136").Value = 3594.75
$ws.Range("L136").Value = 2807.25
$ws.Range("M136").Value = -1044.75
$ws.Range("N136").Value = -7907.25
